# Update column F (dSF) values on Sheet1 to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    9  = -5
    14 = -1
    15 = 1
    16 = 5
    25 = 1
    34 = 3
    38 = -2
    44 = 1
    51 = -9
    52 = -3
    53 = -3
    56 = -1
    58 = 2
    62 = -2
    64 = -3
    66 = -5
    67 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
